$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.664.78'
$ws.Range('E2').Value = '  -1.25%  '
$ws.Range('D3').Value = '2.378.64'
$ws.Range('E3').Value = '  +0.49%  '
$ws.Range('E4').Value = '  -0.22%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '331.27'
$ws.Range('E5').Value = '  +7.14%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '99.54'
$ws.Range('E6').Value = '  -6.02%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.626'
$ws.Range('E9').Value = '  -1.03%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '40.32'
$ws.Range('E10').Value = '  -5.89%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0920'
$ws.Range('E11').Value = '  -2.15%  '
$ws.Range('E12').Value = '  -5.40%  '
$ws.Range('E13').Value = '  -5.44%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.105'
$ws.Range('E14').Value = '  +0.25%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '16.36'
$ws.Range('E15').Value = '  -1.25%  '
$ws.Range('D16').Value = '2.740.47'
$ws.Range('E16').Value = '  +0.38%  '
$ws.Range('D17').Value = '2.366.40'
$ws.Range('E17').Value = '  -0.53%  '
$ws.Range('D18').Value = '42.608.25'
$ws.Range('E18').Value = '  -1.36%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.90'
$ws.Range('E19').Value = '  +6.82%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0000107'
$ws.Range('E20').Value = '  -2.15%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '3.76'
$ws.Range('E21').Value = '  +9.79%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '75.15'
$ws.Range('E22').Value = '  -0.29%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '273.69'
$ws.Range('E23').Value = '  +8.52%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.31'
$ws.Range('E24').Value = '  -8.06%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.70'
$ws.Range('E25').Value = '  +8.00%  '
$ws.Range('E26').Value = '  +0.01%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.46'
$ws.Range('E27').Value = '  -4.69%  '
$ws.Range('E28').Value = '  -0.48%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '23.72'
$ws.Range('E29').Value = '  +4.42%  '
$ws.Range('E30').Value = '  -2.25%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '175.27'
$ws.Range('E31').Value = '  +1.31%  '
$ws.Range('E32').Value = '  -2.19%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0906'
$ws.Range('E33').Value = '  -0.45%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '35.33'
$ws.Range('E34').Value = '  -9.37%  '
$ws.Range('E35').Value = '  +3.94%  '
$ws.Range('E36').Value = '  +2.19%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.60'
$ws.Range('E37').Value = '  -7.41%  '
$ws.Range('E38').Value = '  +6.98%  '
$ws.Range('E39').Value = '  -5.00%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.87'
$ws.Range('E40').Value = '  -4.30%  '
$ws.Range('E41').Value = '  +2.21%  '
$ws.Range('E42').Value = '  +0.32%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '69.53'
$ws.Range('E43').Value = '  -3.85%  '
$ws.Range('E44').Value = '  -1.19%  '
$ws.Range('E45').Value = '  -0.17%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '117.49'
$ws.Range('E46').Value = '  +4.49%  '
$ws.Range('B47').Value = 'BitcoinSV'
$ws.Range('C47').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '90.25'
$ws.Range('E47').Value = '  +31.65%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '11.96'
$ws.Range('E48').Value = '  -3.09%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '5.50'
$ws.Range('E49').Value = '  -3.28%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.17'
$ws.Range('E50').Value = '  -2.51%  '
$ws.Range('D51').Value = '1.599.13'
$ws.Range('E51').Value = '  +6.72%  '
